# Update the exploreCSR AI banner for 2024.
$p = $ppt.ActivePresentation

# --- 1. Refresh the cached "datetimeFigureOut" field text (10/1/22 -> 10/8/23)
#        on the slide master and every slide layout's Date placeholder. ---
$master = $p.SlideMaster

function Update-DatePlaceholder($shapes, $newText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "10/1/22") {
                $tr.Text = $newText
            }
        }
    }
}

Update-DatePlaceholder $master.Shapes "10/8/23"

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes "10/8/23"
}

# --- 2. Update the title textbox on slide 1: reposition/resize and retext. ---
$slide = $p.Slides.Item(1)
$shapes = $slide.Shapes
for ($i = 1; $i -le $shapes.Count; $i++) {
    $shp = $shapes.Item($i)
    if ($shp.Name -eq "TextBox 9") {
        $shp.Left = 48.6216
        $shp.Width = 452.1327

        $tr = $shp.TextFrame.TextRange
        $yearRange = $tr.Find(" 2022-2023:")
        if ($yearRange -ne $null) {
            $yearRange.Text = " 2024:"
        }
        $subtitleRange = $tr.Find("Socially-Responsible AI for Computational Creativity")
        if ($subtitleRange -ne $null) {
            $subtitleRange.Text = "Socially-Responsible Artificial Intelligence"
        }
    }
}
